$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal numbers (e.g. "1.00", "446.53") that
# Excel would otherwise auto-detect as numeric when assigned via .Value. The
# source data stores every Price/Volume cell as text, so mark these specific
# cells as Text before writing them, then drop back to the Normal style so no
# extra formatting is left behind on the cell.
$textPriceCells = @(
    "D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D16", "D17",
    "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29",
    "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40",
    "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50"
)
foreach ($c in $textPriceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '54.019.95'
$ws.Range("E2").Value = '  -11.04%  '
$ws.Range("D3").Value = '2.281.73'
$ws.Range("E3").Value = '  -21.46%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '446.53'
$ws.Range("E5").Value = '  -15.50%  '
$ws.Range("D6").Value = '128.25'
$ws.Range("E6").Value = '  -11.51%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.470'
$ws.Range("E8").Value = '  -15.46%  '
$ws.Range("D9").Value = '2.285.42'
$ws.Range("E9").Value = '  -21.67%  '
$ws.Range("D10").Value = '5.36'
$ws.Range("E10").Value = '  -10.81%  '
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  -15.85%  '
$ws.Range("D12").Value = '0.307'
$ws.Range("E12").Value = '  -15.89%  '
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").Value = '2.678.09'
$ws.Range("E14").Value = '  -21.59%  '
$ws.Range("D15").Value = '54.060.40'
$ws.Range("E15").Value = '  -10.93%  '
$ws.Range("D16").Value = '18.60'
$ws.Range("E16").Value = '  -18.34%  '
$ws.Range("D17").Value = '0.0000119'
$ws.Range("E17").Value = '  -15.94%  '
$ws.Range("D18").Value = '2.300.39'
$ws.Range("E18").Value = '  -21.03%  '
$ws.Range("D19").Value = '4.03'
$ws.Range("E19").Value = '  -20.22%  '
$ws.Range("D20").Value = '299.50'
$ws.Range("E20").Value = '  -17.27%  '
$ws.Range("D21").Value = '9.34'
$ws.Range("E21").Value = '  -20.15%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '5.64'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").Value = '5.26'
$ws.Range("E24").Value = '  -20.60%  '
$ws.Range("D25").Value = '55.44'
$ws.Range("E25").Value = '  -14.43%  '
$ws.Range("D26").Value = '0.977'
$ws.Range("E26").Value = '  -2.13%  '
$ws.Range("D27").Value = '0.158'
$ws.Range("E27").Value = '  -12.71%  '
$ws.Range("D28").Value = '0.367'
$ws.Range("E28").Value = '  -19.48%  '
$ws.Range("B29").Value = 'USDe'
$ws.Range("C29").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '6.79'
$ws.Range("E30").Value = '  -13.67%  '
$ws.Range("D31").Value = '0.0₃0702'
$ws.Range("E31").Value = '  -18.69%  '
$ws.Range("D32").Value = '144.20'
$ws.Range("E32").Value = '  -4.57%  '
$ws.Range("D33").Value = '16.86'
$ws.Range("E33").Value = '  -14.64%  '
$ws.Range("D34").Value = '1.33'
$ws.Range("E34").Value = '  -20.75%  '
$ws.Range("D35").Value = '4.69'
$ws.Range("E35").Value = '  -15.83%  '
$ws.Range("D36").Value = '3.61'
$ws.Range("E36").Value = '  -17.93%  '
$ws.Range("D37").Value = '0.834'
$ws.Range("E37").Value = '  -17.09%  '
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -17.25%  '
$ws.Range("D39").Value = '0.995'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("D40").Value = '32.80'
$ws.Range("E40").Value = '  -13.02%  '
$ws.Range("D41").Value = '10.30'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").Value = '1.24'
$ws.Range("E42").Value = '  -16.93%  '
$ws.Range("D43").Value = '3.13'
$ws.Range("E43").Value = '  -16.12%  '
$ws.Range("D44").Value = '1.914.06'
$ws.Range("E44").Value = '  -16.61%  '
$ws.Range("D45").Value = '0.0493'
$ws.Range("E45").Value = '  -15.22%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.511'
$ws.Range("E46").Value = '  -21.30%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0205'
$ws.Range("E47").Value = '  -13.46%  '
$ws.Range("D48").Value = '0.0811'
$ws.Range("E48").Value = '  -12.31%  '
$ws.Range("D49").Value = '16.13'
$ws.Range("E49").Value = '  -22.01%  '
$ws.Range("D50").Value = '4.01'
$ws.Range("E50").Value = '  -20.55%  '
$ws.Range("E51").Value = '  -3.23%  '

foreach ($c in $textPriceCells) {
    $ws.Range($c).Style = "Normal"
}
